$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving data row (row 2) with the new id and date values
$ws.Range("A2").Value = 253295
$ws.Range("B2").Value = 45911.58333333334

# Remove rows 3 through 16 which are no longer present in the final version
$ws.Range("A3:B16").EntireRow.Delete()
